# ===================================================================
# Offene Fragen.docx - apply edits described by the commit
#   "Checken Sciebo datenbank, hochladen Energiebericht UIW bearbeitet,
#    Latex erste Beschreibung Gebäudemodell"
# ===================================================================

$d = $word.ActiveDocument

# -------------------------------------------------------------------
# 1) Remove the stray "_GoBack" bookmark that used to sit right after
#    the "Offene Fragen" title.
# -------------------------------------------------------------------
$d.Bookmarks("_GoBack").Delete()

# -------------------------------------------------------------------
# 2) "Check en wie groß die Leistung von der Kältemaschine ist,"
#    -> colour everything except the trailing comma green (00B050),
#       leaving the comma itself in a separate, unformatted run.
# -------------------------------------------------------------------
$full = $d.Content
$found = $full.Find.Execute("Check en wie groß die Leistung von der Kältemaschine ist,")
if ($found) {
    $mainRange = $d.Range($full.Start, $full.End - 1)
    $mainRange.Font.Color = 5287936
}

# -------------------------------------------------------------------
# 3) Insert a batch of new bullet points right before the existing
#    "Großer Unterschied zwischen Gebäuden" bullet.
# -------------------------------------------------------------------
function Add-ItemBeforeGrosserUnterschied($texts, $level) {
    $f = $d.Content
    $f.Find.Execute("Großer Unterschied zwischen Gebäuden") | Out-Null
    $anchorPara = $d.Range($f.Start, $f.Start).Paragraphs(1)
    $anchorPara.Range.InsertParagraphBefore()

    $f2 = $d.Content
    $f2.Find.Execute("Großer Unterschied zwischen Gebäuden") | Out-Null
    $anchorPara2 = $d.Range($f2.Start, $f2.Start).Paragraphs(1)
    $newPara = $anchorPara2.Previous()

    if ($level -ne 0) {
        $newPara.Range.ListFormat.ListLevelNumber = $level + 1
    }

    $p = $newPara.Range.Start
    foreach ($t in $texts) {
        $ins = $d.Range($p, $p)
        $ins.InsertAfter($t)
        $p = $p + $t.Length
    }
}

Add-ItemBeforeGrosserUnterschied @("Fotodokumentation HÖRN Zone D  / Begehung Technik ") 0
Add-ItemBeforeGrosserUnterschied @("Sciebo noch angucken : ", [string][char]0x2192, " Datenbank") 0
Add-ItemBeforeGrosserUnterschied @("Bericht von UIW falsch?!") 0
Add-ItemBeforeGrosserUnterschied @("Lüftung auch im Flur? ", [string][char]0x2192, " andere Angaben in Datenbank") 1
Add-ItemBeforeGrosserUnterschied @("würde große Auslegung erklären (statt >60m³/m² eher 30)") 1
Add-ItemBeforeGrosserUnterschied @("Bruttogeschosshöhe 5,1m", " (5,05 Datenbank)") 1
Add-ItemBeforeGrosserUnterschied @("WARUM Wärmerückgewinnung kein Vorteil??") 0

# the "_GoBack" bookmark moved down onto the new "Bruttogeschosshöhe"
# bullet (the new last-edited spot in the document)
$gb = $d.Content
$gb.Find.Execute("Bruttogeschosshöhe 5,1m (5,05 Datenbank)") | Out-Null
$bmRange = $d.Range($gb.End, $gb.End)
$d.Bookmarks.Add("_GoBack", $bmRange) | Out-Null

# -------------------------------------------------------------------
# 4) "An welchen Fernwärmestrang angeschlossen?" -> append
#    " HeißWwasser 1"
# -------------------------------------------------------------------
$full4 = $d.Content
$found4 = $full4.Find.Execute("An welchen Fernwärmestrang angeschlossen?")
if ($found4) {
    $ins4 = $d.Range($full4.End, $full4.End)
    $ins4.InsertAfter(" HeißWwasser 1")
}
